# Applies the Vaccine delivery data refresh described in the commit diff.
# Updates status/date/unit/dose figures across several rows in the Pfizer,
# Moderna, AstraZeneca and J&J blocks, then removes the now-obsolete last
# data row (old row 95, Sputnik) after folding its contents up into row 94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pfizer block ---------------------------------------------------------
# Row 28: status Confirmed -> Delivered
$ws.Cells.Item(28, 4).Value = "Delivered"

# Row 29: status Confirmed -> Delivered
$ws.Cells.Item(29, 4).Value = "Delivered"

# Row 30: Planned -> Delivered, units/doses updated, status-update date moves
$ws.Cells.Item(30, 2).Value = 227
$ws.Cells.Item(30, 3).Value = 265590
$ws.Cells.Item(30, 4).Value = "Delivered"
$ws.Cells.Item(30, 5).Value = 44376

# Row 31: Assumption -> Confirmed, units/doses revised, status-update date set
$ws.Cells.Item(31, 2).Value = 71
$ws.Cells.Item(31, 3).Value = 83070
$ws.Cells.Item(31, 4).Value = "Confirmed"
$ws.Cells.Item(31, 5).Value = 44376

# Row 32: units/doses revised (status untouched)
$ws.Cells.Item(32, 2).Value = 71
$ws.Cells.Item(32, 3).Value = 83070

# Row 33: units/doses revised (status untouched)
$ws.Cells.Item(33, 2).Value = 71
$ws.Cells.Item(33, 3).Value = 83070

# --- Moderna block ---------------------------------------------------------
# Row 49: status Confirmed -> Delivered
$ws.Cells.Item(49, 4).Value = "Delivered"

# Row 50: units revised, status Confirmed -> Delivered
$ws.Cells.Item(50, 2).Value = 384
$ws.Cells.Item(50, 4).Value = "Delivered"

# Row 51: delivery date shifts a day, status Planned -> Delivered
$ws.Cells.Item(51, 1).Value = 44372
$ws.Cells.Item(51, 4).Value = "Delivered"

# Row 52: delivery date shifts a day, status Planned -> Confirmed, status-update date set
$ws.Cells.Item(52, 1).Value = 44379
$ws.Cells.Item(52, 4).Value = "Confirmed"
$ws.Cells.Item(52, 5).Value = 44372

# Row 53: delivery date shifts a day
$ws.Cells.Item(53, 1).Value = 44386

# Row 54: delivery date shifts a day
$ws.Cells.Item(54, 1).Value = 44393

# Row 55: delivery date shifts a day
$ws.Cells.Item(55, 1).Value = 44400

# Row 56: delivery date shifts a day
$ws.Cells.Item(56, 1).Value = 44407

# --- AstraZeneca block -------------------------------------------------------
# Row 75: Assumption -> Delivered, units/doses revised, status-update date set
$ws.Cells.Item(75, 2).Value = 5760
$ws.Cells.Item(75, 3).Value = 57600
$ws.Cells.Item(75, 4).Value = "Delivered"
$ws.Cells.Item(75, 5).Value = 44363

# Row 76: doses revised, Assumption -> Confirmed
$ws.Cells.Item(76, 3).Value = 328800
$ws.Cells.Item(76, 4).Value = "Confirmed"

# --- J&J block ---------------------------------------------------------
# Row 88: status Confirmed -> Delivered, status-update date shifts
$ws.Cells.Item(88, 4).Value = "Delivered"
$ws.Cells.Item(88, 5).Value = 44362

# Row 89: delivery date shifts, doses revised, Assumption -> Delivered, status-update date set
$ws.Cells.Item(89, 1).Value = 44376
$ws.Cells.Item(89, 3).Value = 22615
$ws.Cells.Item(89, 4).Value = "Delivered"
$ws.Cells.Item(89, 5).Value = 44372

# Row 90: delivery date shifts, doses revised, Assumption -> Confirmed
$ws.Cells.Item(90, 1).Value = 44383
$ws.Cells.Item(90, 3).Value = 16800
$ws.Cells.Item(90, 4).Value = "Confirmed"

# Row 91: delivery date shifts, doses revised
$ws.Cells.Item(91, 1).Value = 44390
$ws.Cells.Item(91, 3).Value = 23000

# Row 92: delivery date shifts, doses revised
$ws.Cells.Item(92, 1).Value = 44397
$ws.Cells.Item(92, 3).Value = 23000

# Row 93: delivery date shifts, doses revised
$ws.Cells.Item(93, 1).Value = 44404
$ws.Cells.Item(93, 3).Value = 23000

# Row 94 previously held a J&J Assumption row; the source data dropped it and
# moved the trailing Sputnik row (old row 95) up into row 94's place.
# (Column E is already blank text in this row, so it's left untouched.)
$ws.Cells.Item(94, 1).Value = 44256
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(94, 3).Value = 200000
$ws.Cells.Item(94, 4).Value = "Delivered"
$ws.Cells.Item(94, 6).Value = "Sputnik"

# Remove the now-duplicate trailing row so the sheet ends at row 94.
$ws.Rows.Item(95).Delete()
